# Consolidate text-run nodes on slide 1:
#   - Title shape: "A" " " "Table," " " "with" " " "a" " " "caption"
#       -> "A " "Table, " "with " "a " "caption"
#   - Caption textbox: "Demonstration" " " "of" " " "simple" " " "table" " "
#     "syntax," " " "with" " " "alignment"
#       -> "Demonstration " "of " "simple " "table " "syntax, " "with " "alignment"
#
# Each merge is done by re-writing a Characters() sub-range so the
# underlying run model collapses the word + following space into a
# single run, instead of leaving them as separate runs.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Title shape: "A Table, with a caption" ---
$title = $s.Shapes.Item(1).TextFrame.TextRange
$title.Characters(1, 2).Text = "A "
$title.Characters(3, 7).Text = "Table, "
$title.Characters(10, 5).Text = "with "
$title.Characters(15, 2).Text = "a "

# --- Caption textbox: "Demonstration of simple table syntax, with alignment" ---
$caption = $s.Shapes.Item(3).TextFrame.TextRange
$caption.Characters(1, 14).Text = "Demonstration "
$caption.Characters(15, 3).Text = "of "
$caption.Characters(18, 7).Text = "simple "
$caption.Characters(25, 6).Text = "table "
$caption.Characters(31, 8).Text = "syntax, "
$caption.Characters(39, 5).Text = "with "
